$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B176").Value = 57552
$ws.Range("E176").Value = 136.86
$ws.Range("F176").Value = -5
$ws.Range("G176").Value = -603.45
$ws.Range("B177").Value = 64329
$ws.Range("E177").Value = 128.32
$ws.Range("F177").Value = 6
$ws.Range("G177").Value = 724.14
$ws.Range("B256").Value = 48719
$ws.Range("E256").Value = 353.35
$ws.Range("F256").Value = -81
$ws.Range("G256").Value = -23955.75
$ws.Range("B257").Value = 64979
$ws.Range("E257").Value = 314.41
$ws.Range("F257").Value = 82
$ws.Range("G257").Value = 24251.5
$ws.Range("B271").Value = 48706
$ws.Range("E271").Value = 39.8
$ws.Range("F271").Value = -144
$ws.Range("G271").Value = -4795.2
$ws.Range("B272").Value = 64973
$ws.Range("E272").Value = 35.4
$ws.Range("F272").Value = 150
$ws.Range("G272").Value = 4995
$ws.Range("B305").Value = 57854
$ws.Range("F305").Value = 2
$ws.Range("G305").Value = 611.6799999999999
$ws.Range("B306").Value = 62997
$ws.Range("F306").Value = 72
$ws.Range("G306").Value = 22020.48
$ws.Range("B338").Value = 55373
$ws.Range("E338").Value = 163.62
$ws.Range("F338").Value = -94
$ws.Range("G338").Value = -13562.32
$ws.Range("B339").Value = 63520
$ws.Range("E339").Value = 153.4
$ws.Range("F339").Value = 97
$ws.Range("G339").Value = 13995.16
$ws.Range("B342").Value = 57802
$ws.Range("E342").Value = 162.71
$ws.Range("F342").Value = -79
$ws.Range("G342").Value = -11334.92
$ws.Range("B343").Value = 63571
$ws.Range("F343").Value = 29
$ws.Range("G343").Value = 4160.92
$ws.Range("B344").Value = 63531
$ws.Range("E344").Value = 152.53
$ws.Range("F344").Value = 80
$ws.Range("G344").Value = 11478.4
$ws.Range("B364").Value = 57885
$ws.Range("E364").Value = 62.28
$ws.Range("F364").Value = 4
$ws.Range("G364").Value = 208.52
$ws.Range("B365").Value = 63652
$ws.Range("E365").Value = 55.42
$ws.Range("F365").Value = 250
$ws.Range("G365").Value = 13032.5
$ws.Range("B392").Value = 57835
$ws.Range("F392").Value = 1
$ws.Range("G392").Value = 59.13
$ws.Range("B393").Value = 62933
$ws.Range("F393").Value = 146
$ws.Range("G393").Value = 8632.98
$ws.Range("B413").Value = 57857
$ws.Range("F413").Value = 3
$ws.Range("G413").Value = 453.51
$ws.Range("B414").Value = 63008
$ws.Range("F414").Value = 504
$ws.Range("G414").Value = 76189.67999999999
$ws.Range("B423").Value = 63102
$ws.Range("C423").Value = "HUL-Vim Bar Multipack Fw 4X200G"
$ws.Range("F423").Value = 36
$ws.Range("G423").Value = 2140.92
$ws.Range("B424").Value = 53082
$ws.Range("C424").Value = "HUL-VIM BAR MULTIPACK FW 4X200G"
$ws.Range("F424").Value = 1
$ws.Range("G424").Value = 59.47
$ws.Range("B528").Value = 47097
$ws.Range("D528").Value = 112.28
$ws.Range("E528").Value = 134.16
$ws.Range("F528").Value = 15
$ws.Range("G528").Value = 1684.2
$ws.Range("B529").Value = 58047
$ws.Range("D529").Value = 105.54
$ws.Range("E529").Value = 126.1
$ws.Range("F529").Value = 54
$ws.Range("G529").Value = 5699.16
$ws.Range("B578").Value = 64915
$ws.Range("E578").Value = 20.98
$ws.Range("F578").Value = 40
$ws.Range("G578").Value = 789.2
$ws.Range("B579").Value = 45695
$ws.Range("E579").Value = 23.58
$ws.Range("F579").Value = -36
$ws.Range("G579").Value = -710.28
$ws.Range("B585").Value = 45718
$ws.Range("E585").Value = 19.38
$ws.Range("F585").Value = -294
$ws.Range("G585").Value = -4768.68
$ws.Range("B586").Value = 64927
$ws.Range("E586").Value = 17.26
$ws.Range("F586").Value = 295
$ws.Range("G586").Value = 4784.9
$ws.Range("B593").Value = 64919
$ws.Range("E593").Value = 27.97
$ws.Range("F593").Value = 224
$ws.Range("G593").Value = 5891.2
$ws.Range("B594").Value = 45702
$ws.Range("E594").Value = 31.43
$ws.Range("F594").Value = -215
$ws.Range("G594").Value = -5654.5
$ws.Range("B712").Value = 60022
$ws.Range("E712").Value = 37.22
$ws.Range("F712").Value = -113
$ws.Range("G712").Value = -3709.79
$ws.Range("B713").Value = 64830
$ws.Range("E713").Value = 34.9
$ws.Range("F713").Value = 117
$ws.Range("G713").Value = 3841.11
$ws.Range("B864").Value = 54751
$ws.Range("E864").Value = 46.34
$ws.Range("F864").Value = -19
$ws.Range("G864").Value = -776.53
$ws.Range("B865").Value = 65079
$ws.Range("E865").Value = 43.44
$ws.Range("F865").Value = 21
$ws.Range("G865").Value = 858.27

Write-Host "Applied stock report swap updates"
